$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text formatting (avoid numeric auto-conversion)
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D8', 'D10', 'D11', 'D12', 'D13', 'D14', 'D16', 'D17', 'D18', 'D19', 'D20', 'D22', 'D23', 'D25', 'D26', 'D28', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D40', 'D41', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.350.16'
$ws.Range('E2').Value = '  +3.98%  '
$ws.Range('D3').Value = '1.593.30'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '214.42'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = '0.495'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').Value = '24.06'
$ws.Range('E8').Value = '  +8.08%  '
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').Value = '0.0600'
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('D11').Value = '0.0887'
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('D12').Value = '1.822.68'
$ws.Range('E12').Value = '  +1.87%  '
$ws.Range('D13').Value = '1.593.08'
$ws.Range('E13').Value = '  +2.21%  '
$ws.Range('D14').Value = '0.531'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '28.374.22'
$ws.Range('E16').Value = '  +4.19%  '
$ws.Range('D17').Value = '63.12'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('D18').Value = '227.18'
$ws.Range('E18').Value = '  +4.52%  '
$ws.Range('D19').Value = '0.0₃0710'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('D20').Value = '7.49'
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = '4.09'
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('D23').Value = '9.32'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('D25').Value = '151.74'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').Value = '15.20'
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').Value = '6.58'
$ws.Range('E28').Value = '  -0.75%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +0.74%  '
$ws.Range('D31').Value = '0.0475'
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').Value = '3.23'
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('D33').Value = '3.13'
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').Value = '1.394.19'
$ws.Range('E34').Value = '  -4.38%  '
$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('D36').Value = '1.04'
$ws.Range('E36').Value = '  -5.82%  '
$ws.Range('E37').Value = '  +0.61%  '
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('E39').Value = '  +7.37%  '
$ws.Range('D40').Value = '0.539'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('D41').Value = '0.814'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  +6.81%  '
$ws.Range('D45').Value = '0.985'
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').Value = '64.41'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').Value = '1.733.22'
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '87.45'
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('B49').Value = 'mCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D49').Value = '2.14'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').Value = '0.0₆0104'
$ws.Range('E50').Value = '  +3.84%  '
$ws.Range('D51').Value = '0.0524'
$ws.Range('E51').Value = '  -0.12%  '
